# Update the Gantt chart's "Display Week" pointer and refresh a handful
# of task progress percentages on the ProjectSchedule sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# Display Week: jump the visible 8-week window back to week 1.
$ws.Range("E4").Value = 1

# Task progress updates (PROGRESS column, D).
$ws.Range("D9").Value  = 1      # Reflechir sur le design et le fonctionnement du robot -> 100%
$ws.Range("D10").Value = 1      # Modelisation 3D du cannon du robot -> 100%
$ws.Range("D11").Value = 0.5    # Modelisation 3D de la magazine -> 50%
$ws.Range("D12").Value = 1      # Modelisation 3D de la base de la tourelle -> 100%
$ws.Range("D18").Value = 0.5    # controle avec l'application? -> 50%

# Leave the selection on D13, matching where the editor ended up.
$ws.Range("D13").Select()
